# Daily attendance processing - 2026-01-25 10:01:12
# Reorder the "Recorded By" (column G) values for the rows that list
# a flagged / non-"System" account first, so that "System" (or "system")
# is listed first instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Rows (in column G) whose "Recorded By" value needs to be reordered.
$rows = @(2, 3, 4, 5, 6, 7, 8, 28, 29, 30, 31, 32, 33, 34, 54, 55, 56, 57, 58, 59, 60, 80, 81, 82, 106, 107, 108, 132, 133, 134)

foreach ($row in $rows) {
    $cell = $ws.Cells.Item($row, 7)   # column G
    $current = $cell.Value()
    if ($current) {
        $parts = $current -split ",\s*"
        $reversed = $parts[($parts.Length - 1)..0]
        $cell.Value = [string]::Join(", ", $reversed)
    }
}
